$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") for rows 2-261 is being updated
# from serial date 45186 (2023-09-17) to 45188 (2023-09-19).
$newDate = 45188

for ($row = 2; $row -le 261; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
